$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 22:50"

# Update Estados Unidos (row 4)
$ws.Range("B4").Value = 139217
$ws.Range("C4").Value = 15639
$ws.Range("D4").Value = 4435
$ws.Range("E4").Value = 132342
$ws.Range("G4").Value = 220
$ws.Range("H4").Value = 2440

# Update España (row 7)
$ws.Range("B7").Value = 80031
$ws.Range("C7").Value = 6796
$ws.Range("E7").Value = 58520
$ws.Range("G7").Value = 820
$ws.Range("H7").Value = 6802

# Update Alemania (row 8)
$ws.Range("B8").Value = 62095
$ws.Range("C8").Value = 4400
$ws.Range("E8").Value = 52359
$ws.Range("G8").Value = 92
$ws.Range("H8").Value = 525

# Update Marruecos (row 66)
$ws.Range("B66").Value = 479
$ws.Range("C66").Value = 77
$ws.Range("E66").Value = 440

# "Costa de Marfil" moves up in the ranking (inserted right after Oman, row 93),
# so it now occupies row 94 with updated stats, and the countries that used to sit
# between Oman and Costa de Marfil (Islas Feroe, Ghana, Malta, Uzbekistan, Senegal)
# each shift down by one row, keeping their own (unchanged) stats.

# Row 94: Costa de Marfil (new stats, new position)
$ws.Range("A94").Value = "Costa de Marfil"
$ws.Range("B94").Value = 165
$ws.Range("C94").Value = 64
$ws.Range("D94").Value = 4
$ws.Range("E94").Value = 160
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 1

# Row 95: Islas Feroe (shifted down from row 94)
$ws.Range("A95").Value = "Islas Feroe"
$ws.Range("B95").Value = 159
$ws.Range("C95").Value = 4
$ws.Range("D95").Value = 70
$ws.Range("E95").Value = 89
$ws.Range("F95").Value = 1
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0

# Row 96: Ghana (shifted down from row 95)
$ws.Range("A96").Value = "Ghana"
$ws.Range("B96").Value = 152
$ws.Range("C96").Value = 11
$ws.Range("D96").Value = 2
$ws.Range("E96").Value = 145
$ws.Range("F96").Value = 1
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 5

# Row 97: Malta (shifted down from row 96)
$ws.Range("A97").Value = "Malta"
$ws.Range("B97").Value = 151
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 2
$ws.Range("E97").Value = 149
$ws.Range("F97").Value = 4
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0

# Row 98: Uzbekistan (shifted down from row 97)
$ws.Range("A98").Value = "Uzbekistan"
$ws.Range("B98").Value = 144
$ws.Range("C98").Value = 40
$ws.Range("D98").Value = 7
$ws.Range("E98").Value = 135
$ws.Range("F98").Value = 8
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 2

# Row 99: Senegal (shifted down from row 98)
$ws.Range("A99").Value = "Senegal"
$ws.Range("B99").Value = 142
$ws.Range("C99").Value = 12
$ws.Range("D99").Value = 27
$ws.Range("E99").Value = 115
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
